$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated vm_pu results for the "case with 380 kV" run.
# Columns mapped to Cells.Item column index: B=2 C=3 D=4 E=5 F=6 I=9 J=10 K=11 L=12 M=13
$colIdx = @(2,3,4,5,6,9,10,11,12,13)

$rowData = @{
    2 = @(1.02,1.048995049297711,1.055427139596091,1.05617780124415,1.066690171761095,1.049910084595802,1.054035845434294,1.058167627690303,1.058916226787732,1.069400043080337)
    3 = @(1.02,1.049845166754432,1.05608685637099,1.056915886730985,1.067465853250328,1.050145554399191,1.054535551589984,1.058641175864717,1.059468092199708,1.069991465215181)
    4 = @(1.02,1.050395919622084,1.056514287087109,1.057394412396004,1.067968690458155,1.050296967654735,1.05485888137687,1.058947456738365,1.059825450192252,1.07037441003271)
    5 = @(1.02,1.050627614821044,1.056694108769313,1.057595806529154,1.068180301347833,1.050360393150235,1.054994804642434,1.059076183188099,1.059975745253873,1.070535459514143)
    6 = @(1.02,1.050666526720588,1.056724309179956,1.057629634449963,1.068215844480775,1.050371029137728,1.055017626449397,1.059097794890362,1.060000984034707,1.070562503872052)
    7 = @(1.02,1.050399014923256,1.056516689365462,1.057397102564567,1.067971517158619,1.050297816048767,1.054860697609863,1.058949176922905,1.059827458201399,1.070376561752059)
    8 = @(1.02,1.049282210708367,1.055649979021987,1.056427045951694,1.066952125389382,1.04998985920692,1.054204725459441,1.058327692950028,1.059102676699379,1.069599863165911)
    9 = @(1.02,1.047319466420205,1.054127016535561,1.054724927129504,1.065162956127726,1.049439950881012,1.053048774617573,1.057231573827143,1.057827607789804,1.068233246038715)
    10 = @(1.02,1.046014573491201,1.053114702773835,1.05359516556281,1.063975088296711,1.049068522690467,1.052278191139925,1.056500246684689,1.056979054540224,1.067323625121013)
    11 = @(1.02,1.045450416575334,1.052677093046873,1.053107172505134,1.063461918425263,1.048906557100461,1.051944548836872,1.05618345327096,1.056611994326801,1.066930115050192)
    12 = @(1.02,1.045240996060806,1.052514656673263,1.052926092562685,1.063271484311816,1.048846226241561,1.051820624408913,1.056065764850078,1.056475708757539,1.066784003875913)
    13 = @(1.02,1.045285911449193,1.052549494757407,1.052964926539489,1.063312324896443,1.048859175077862,1.051847206363215,1.056091010173964,1.056504939874041,1.06681534265905)
    14 = @(1.02,1.045433103097669,1.052663663722982,1.053092200648303,1.063446173397824,1.048901573588245,1.051934305091495,1.056173725453136,1.056600727748909,1.066918036303391)
    15 = @(1.02,1.045523810337195,1.052734021803304,1.053170642656223,1.063528665792989,1.048927674258124,1.051987970219492,1.056224686823092,1.056659753409854,1.066981316714628)
    16 = @(1.02,1.046052033211828,1.053143761004901,1.053627577528097,1.064009170825574,1.049079247978556,1.052300334532454,1.05652126872669,1.057003423019638,1.067349748850254)
    17 = @(1.02,1.046383607695574,1.053400976023854,1.053914523659798,1.064310897560919,1.049174022970679,1.052496280180197,1.056707274511174,1.057219097727448,1.067580954608863)
    18 = @(1.02,1.046577093340819,1.053551075387345,1.054082010247354,1.064487003834695,1.049229194077556,1.05261057432369,1.056815756487958,1.05734493263218,1.067715847814144)
    19 = @(1.02,1.046643081117173,1.053602267240757,1.054139138421736,1.064547070827075,1.049247987392995,1.052649546032353,1.056852743952743,1.057387845057785,1.067761848763861)
    20 = @(1.02,1.046348024205112,1.053373372012539,1.053883725067471,1.064278513313532,1.049163865835312,1.052475256812459,1.056687319087321,1.057195954189967,1.067556144807979)
    21 = @(1.02,1.045389755154793,1.05263004071904,1.053054716544527,1.063406753364352,1.048889092959372,1.051908656533024,1.056149368329936,1.056572519023355,1.066887794026083)
    22 = @(1.02,1.044788020732124,1.05216332407012,1.052534542296138,1.062859685862046,1.048715351359386,1.051552443084685,1.055811038263825,1.056180870757871,1.066467900094046)
    23 = @(1.02,1.045106938220667,1.052410677612114,1.052810195731374,1.063149597193443,1.048807547769623,1.051741275215844,1.055990402364243,1.056388459114254,1.066690462486442)
    24 = @(1.02,1.046364102593877,1.053385844860404,1.053897641268197,1.064293146008448,1.049168455746147,1.05248475636027,1.056696336121795,1.057206411648866,1.067567355181761)
    25 = @(1.02,1.047826255251105,1.054520219000314,1.055164095411742,1.065624642210949,1.049582969142593,1.053347612594776,1.057515055252496,1.058156987036892,1.068586299489538)
}

foreach ($row in $rowData.Keys) {
    $values = $rowData[$row]
    for ($i = 0; $i -lt $colIdx.Count; $i++) {
        $ws.Cells.Item($row, $colIdx[$i]).Value = $values[$i]
    }
}

Write-Output "Updated $($rowData.Keys.Count) rows"

